# Weekly price-sheet update: a new daily observation (date serial 45205)
# is inserted as a new row at position 57; all existing rows from 57
# downward shift down by one (row 90 -> row 91), and the new row 57
# carries the same market/volume/price data as the row it was copied
# from (the row immediately below it once shifted), only the date
# (column D) differs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 57; everything from old row 57 onward moves to 58..91.
$ws.Rows.Item(57).Insert()

# Seed the new row 57 with the same data as the (now shifted) row 58,
# then overwrite just the date so the rest of the fields line up.
$ws.Range("A58:R58").Copy()
$ws.Range("A57:R57").PasteSpecial()

$ws.Range("D57").Value = 45205
